$d = $word.ActiveDocument
$table = $d.Tables(1)

# Map of "row,col" -> new text, derived from the diff (row-major order,
# skipping the blank spacer rows). The table has duplicate "92÷9=" cells
# so we must address them by position rather than by global Find/Replace.
$updates = @{
    "1,1" = "66÷4="
    "1,2" = "88÷2="
    "1,3" = "14÷7="
    "1,4" = "49÷8="
    "1,5" = "59÷8="
    "5,1" = "25÷9="
    "5,2" = "86÷4="
    "5,3" = "65÷7="
    "5,4" = "50÷6="
    "5,5" = "98÷2="
    "9,1" = "23÷8="
    "9,2" = "35÷2="
    "9,3" = "27÷8="
    "9,4" = "99÷2="
    "9,5" = "30÷9="
    "13,1" = "20÷3="
    "13,2" = "49÷9="
    "13,3" = "14÷8="
    "13,4" = "99÷3="
    "13,5" = "31÷2="
    "17,1" = "43÷3="
    "17,2" = "55÷7="
    "17,3" = "12÷8="
    "17,4" = "52÷2="
    "17,5" = "21÷7="
}

foreach ($key in $updates.Keys) {
    $parts = $key.Split(",")
    $r = [int]$parts[0]
    $c = [int]$parts[1]
    $cell = $table.Cell($r, $c)
    $rng = $cell.Range
    $rng.End = $rng.End - 1
    $rng.Text = $updates[$key]
}
